# Weekly fruit/vegetable price update: a new weekly record was inserted
# into the daily-logic subset sheet for "Ají" (Agrícola del Norte S.A. de
# Arica), pushing the existing rows 78-100 down to 79-101 and adding a
# brand-new row 78 with this week's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 78 - this shifts rows 78:100
# down to 79:101 (and copies row 77's formatting, e.g. the date style on
# column D, down onto the new row), matching the dimension growing from
# A1:R100 to A1:R101.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new weekly record.
$ws.Range("A78").Value = 1
$ws.Range("B78").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C78").Value = "Arica y Parinacota"
$ws.Range("D78").Value = 44809
$ws.Range("E78").Value = 15
$ws.Range("F78").Value = 100112021
$ws.Range("G78").Value = "Ají"
$ws.Range("H78").Value = "Inferno"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 130
$ws.Range("K78").Value = 14000
$ws.Range("L78").Value = 15000
$ws.Range("M78").Value = 14500
$ws.Range("N78").Value = "$/caja 15 kilos"
$ws.Range("O78").Value = "Región de Arica y Parinacota"
$ws.Range("P78").Value = 967
$ws.Range("Q78").Value = 15
$ws.Range("R78").Value = "Hortaliza"
